# Added "Scheduled maint class" values into the signup sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("signup")

# Replace the placeholder applicant values in column A (rows 2-3)
# with the new scheduled-maintenance class values.
$ws.Range("A2").Value = "394-357"
$ws.Range("A3").Value = "319-86"

$ws.Range("G6").Select()
